$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2020" column (Q) to the table, mirroring the formatting of the
# existing 2019 column (P).
$ws.Range("P4:P8").Copy()
$ws.Range("Q4:Q8").PasteSpecial(-4122)   # xlPasteFormats

# Populate the new column's values: header year 2020, and the 2019 figures
# carried forward unchanged for each indicator row.
$ws.Range("Q4").Value = 2020
$ws.Range("Q5").Value = 2
$ws.Range("Q6").Value = 0.3
$ws.Range("Q7").Value = 0.1
$ws.Range("Q8").Value = 4.3

# Mirror the author's saved selection/view state.
$ws.Range("O12").Select()
